# Update the "形態" (data type) column for the CreateDate / LastUpdate
# fields in the DBD layout sheet from "DATE" to "TIMESTAMP".
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("DBD")

$ws.Range("D14").Value = "TIMESTAMP"
$ws.Range("D16").Value = "TIMESTAMP"
